$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-AuthorParagraph($paragraph, [string]$text) {
    $xml = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p>'
    [void]$paragraph.Range.InsertXML($xml)
}

# 1) "Adam Sparks" -> "Prepared for: collab_partner (collab_partner@email.com.au)"
$p1 = $d.Paragraphs(3)
if ($p1.Range.Text.TrimEnd([char]13,[char]7) -ne "Adam Sparks") {
    throw "Unexpected paragraph 3 text: $($p1.Range.Text)"
}
Set-AuthorParagraph $p1 "Prepared for: collab_partner (collab_partner@email.com.au)"

# 2) "Email: cbada@curtin.edu.au" -> "Prepared by: Your.Name"
$p2 = $d.Paragraphs(4)
if ($p2.Range.Text.TrimEnd([char]13,[char]7) -ne "Email: cbada@curtin.edu.au") {
    throw "Unexpected paragraph 4 text: $($p2.Range.Text)"
}
Set-AuthorParagraph $p2 "Prepared by: Your.Name"

# 3) "Project Leads: ..." -> "Project Lead: ..." (singular form)
$p3 = $d.Paragraphs(5)
$expected3 = "Project Leads: Curtin University " + [char]0x2013 + " Prof Mark Gibberd, Dr Julia Easton, Prof Adam Sparks"
if ($p3.Range.Text.TrimEnd([char]13,[char]7) -ne $expected3) {
    throw "Unexpected paragraph 5 text: $($p3.Range.Text)"
}
$newLead = "Project Lead: Curtin University " + [char]0x2013 + " Prof Mark Gibberd, Dr Julia Easton, Prof Adam Sparks"
Set-AuthorParagraph $p3 $newLead

# 4) Insert a brand-new "Author"-styled paragraph right after it: "email: cbada@curtin.edu.au"
[void]$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(6)
Set-AuthorParagraph $p4 "email: cbada@curtin.edu.au"

# 5) Update the "Author" paragraph style: the text is no longer bold
$authorStyle = $d.Styles("Author")
$authorStyle.Font.Bold = $false
